$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.36219633333334
$ws.Range("H2").Value = 37.086589
$ws.Range("I2").Value = 0.8692805094072583
$ws.Range("J2").Value = 0.8692805094072584
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 45.95651366666667
$ws.Range("N2").Value = 137.869541
$ws.Range("O2").Value = 0.6189188856627118
$ws.Range("P2").Value = 0.6189188856627118
$ws.Range("Q2").Value = 568.12344474285
$ws.Range("R2").Value = 5113.11100268565
$ws.Range("S2").Value = 0.5380141242106548
$ws.Range("T2").Value = 0.5380141242106549
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.36219633333334
$ws.Range("H3").Value = 37.086589
$ws.Range("I3").Value = 0.8692805094072583
$ws.Range("J3").Value = 0.8692805094072584
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("N3").Value = 20.549744
$ws.Range("O3").Value = 0.09225115688993263
$ws.Range("P3").Value = 0.09225115688993261
$ws.Range("Q3").Value = 84.6799899759129
$ws.Range("R3").Value = 762.1199097832161
$ws.Range("S3").Value = 0.08019213265468954
$ws.Range("T3").Value = 0.08019213265468954
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.36219633333334
$ws.Range("H4").Value = 37.086589
$ws.Range("I4").Value = 0.8692805094072583
$ws.Range("J4").Value = 0.8692805094072584
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 21.446458
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2888299574473556
$ws.Range("P4").Value = 0.2888299574473556
$ws.Range("Q4").Value = 265.1253244505873
$ws.Range("R4").Value = 2386.127920055286
$ws.Range("S4").Value = 0.251074252541914
$ws.Range("T4").Value = 0.251074252541914
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.434409
$ws.Range("H5").Value = 4.303227
$ws.Range("I5").Value = 0.1008642600875229
$ws.Range("J5").Value = 0.1008642600875229
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 45.95651366666667
$ws.Range("N5").Value = 137.869541
$ws.Range("O5").Value = 0.6189188856627118
$ws.Range("P5").Value = 0.6189188856627118
$ws.Range("Q5").Value = 65.92043681208966
$ws.Range("R5").Value = 593.2839313088069
$ws.Range("S5").Value = 0.06242679545656363
$ws.Range("T5").Value = 0.06242679545656364
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.434409
$ws.Range("H6").Value = 4.303227
$ws.Range("I6").Value = 0.1008642600875229
$ws.Range("J6").Value = 0.1008642600875229
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("N6").Value = 20.549744
$ws.Range("O6").Value = 0.09225115688993263
$ws.Range("P6").Value = 0.09225115688993261
$ws.Range("Q6").Value = 9.825579247098666
$ws.Range("R6").Value = 88.43021322388799
$ws.Range("S6").Value = 0.009304844681921047
$ws.Range("T6").Value = 0.009304844681921047
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.434409
$ws.Range("H7").Value = 4.303227
$ws.Range("I7").Value = 0.1008642600875229
$ws.Range("J7").Value = 0.1008642600875229
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.446458
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2888299574473556
$ws.Range("P7").Value = 0.2888299574473556
$ws.Range("Q7").Value = 30.76299237332199
$ws.Range("R7").Value = 276.8669313598979
$ws.Range("S7").Value = 0.02913261994903826
$ws.Range("T7").Value = 0.02913261994903826
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4245766666666667
$ws.Range("H8").Value = 1.27373
$ws.Range("I8").Value = 0.02985523050521867
$ws.Range("J8").Value = 0.02985523050521867
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.95651366666667
$ws.Range("N8").Value = 137.869541
$ws.Range("O8").Value = 0.6189188856627118
$ws.Range("P8").Value = 0.6189188856627118
$ws.Range("Q8").Value = 19.51206338421444
$ws.Range("R8").Value = 175.60857045793
$ws.Range("S8").Value = 0.01847796599549334
$ws.Range("T8").Value = 0.01847796599549334
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4245766666666667
$ws.Range("H9").Value = 1.27373
$ws.Range("I9").Value = 0.02985523050521867
$ws.Range("J9").Value = 0.02985523050521867
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("N9").Value = 20.549744
$ws.Range("O9").Value = 0.09225115688993263
$ws.Range("P9").Value = 0.09225115688993261
$ws.Range("Q9").Value = 2.908313936124444
$ws.Range("R9").Value = 26.17482542512
$ws.Range("S9").Value = 0.00275417955332203
$ws.Range("T9").Value = 0.00275417955332203
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4245766666666667
$ws.Range("H10").Value = 1.27373
$ws.Range("I10").Value = 0.02985523050521867
$ws.Range("J10").Value = 0.02985523050521867
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.446458
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2888299574473556
$ws.Range("P10").Value = 0.2888299574473556
$ws.Range("Q10").Value = 9.105665649446665
$ws.Range("R10").Value = 81.95099084501999
$ws.Range("S10").Value = 0.0086230849564033
$ws.Range("T10").Value = 0.0086230849564033
